# linkedin_events.xlsx — refresh event listing with new LinkedIn/AI events,
# add Meeting URL hyperlinks + attendee counts for two rows, and append a
# new trailing event row. ("fixes and error handling")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "The future of AI Translation: LLMs or Neural MT? feat. João Graça"
$ws.Range("B2").Value = "Wed, Jun 26, 2024, 7:00 PM - 8:00 PM (your local time)"
$ws.Range("C2").Value = "Nimdzi Insights"

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = "AI Insights: Part 3 - LLMs: Considerations When Choosing an LLM"
$ws.Range("B3").Value = "Wed, Jun 19, 2024, 12:00 AM (your local time)"
$ws.Range("C3").Value = "Domo"

# --- Row 4 (gains a Meeting URL hyperlink + attendee count) --------------
$ws.Range("A4").Value = "Unlocking Enterprise Potential: Harnessing Open Source LLMs for Production"
$ws.Range("B4").Value = "Thu, Jun 20, 2024, 8:30 PM - 9:30 PM (your local time)"
$ws.Range("C4").Value = "WalkingTree Technologies"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://us02web.zoom.us/webinar/register/WN_h5C6V61vSruKKCfWd8_nIQ")
$ws.Range("E4").Value = "220 attendees"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = "AI (LLMs) and XR Powered Digital Twins in Healthcare"
$ws.Range("B5").Value = "Wed, Jun 19, 2024, 9:00 PM - 10:00 PM (your local time)"
$ws.Range("C5").Value = "Alex G. Lee, Ph.D. Esq."

# --- Row 6 (old attendee count removed) -----------------------------------
$ws.Range("A6").Value = "Mastering Data Science: The Impact of LLMs"
$ws.Range("B6").Value = "Wed, Jul 17, 2024, 10:00 PM (your local time)"
$ws.Range("C6").Value = "Data Science Dojo"
$ws.Range("E6").ClearContents()

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "AI Governance: Data minimisation & anonymisation while leveraging LLMs"
$ws.Range("B7").Value = "Thu, Jun 20, 2024, 6:00 PM - 6:30 PM (your local time)"
$ws.Range("C7").Value = "TrustWorks"

# --- Row 8 (gains a Meeting URL hyperlink + new attendee count) ----------
$ws.Range("A8").Value = "How to align LLMs to enterprise objectives/policies"
$ws.Range("B8").Value = "Jun 18, 2024, 10:00 PM - Jun 19, 2024, 11:00 AM (your local time)"
$ws.Range("C8").Value = "Snorkel AI"
$ws.Hyperlinks.Add($ws.Range("D8"), "https://snorkel.ai/event/how-to-align-llms-to-enterprise-objectives-policies/")
$ws.Range("E8").Value = "92 attendees"

# --- Row 9 (old attendee count removed) -----------------------------------
$ws.Range("A9").Value = "Pre-Training LLMs on Personal Computers"
$ws.Range("B9").Value = "Wed, Aug 28, 2024, 9:00 PM (your local time)"
$ws.Range("C9").Value = "Data Science Dojo"
$ws.Range("E9").ClearContents()

# --- Row 10 (old attendee count removed) ----------------------------------
$ws.Range("A10").Value = "Building trustworthiness in AI using RAG"
$ws.Range("B10").Value = "Wed, Jun 26, 2024, 9:30 AM (your local time)"
$ws.Range("C10").Value = "Katonic AI"
$ws.Range("E10").ClearContents()

# --- Row 11 (new trailing row) --------------------------------------------
$ws.Range("A11").Value = "Hybrid Chatbots: Merging LLMs and Classification Models"
$ws.Range("B11").Value = "Wed, Aug 21, 2024, 10:00 PM (your local time)"
$ws.Range("C11").Value = "Data Science Dojo"
